# Insert a new row at 130, shifting existing rows 130-243 down to 131-244.
# The new row 130 is a copy of the (old) row 130 data, but with an updated
# date (column D) and volume (column M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("130:130").Insert()

$row = 130
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44512
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 125
$ws.Cells.Item($row, 14).Value = 8000
$ws.Cells.Item($row, 15).Value = 8000
$ws.Cells.Item($row, 16).Value = 8000
$ws.Cells.Item($row, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item($row, 18).Value = "Perú"
$ws.Cells.Item($row, 19).Value = 2000
$ws.Cells.Item($row, 20).Value = 4
